$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.943.32"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.638.28"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'214.58"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'0.5067"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "'0.2570"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.06352"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "'19.83"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.642.25"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'0.5465"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "0.0₅7735"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "'64.19"
$ws.Range("D17").Value = "25.963.66"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'4.460"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "'195.90"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'9.952"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "'6.140"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "'1.897"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "'142.85"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'0.1256"
$ws.Range("E26").Value = "  +9.82%  "
$ws.Range("D27").Value = "'6.846"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").Value = "'15.64"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'1.236"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "'0.04882"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").Value = "'3.252"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'3.202"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "'1.551"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "'2.376"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "'0.9158"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").Value = "'2.570"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "1.135.72"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'0.5527"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").Value = "'5.589"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "'98.58"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -8.94%  "
$ws.Range("D45").Value = "1.774.02"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "'0.4505"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "'55.22"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "'0.05181"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").Value = "'7.520"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("E51").Value = "  -0.49%  "
